# Trade #160 closed at 2026-02-18 00:44:17 - unknown UNKNOWN +0.000%
#
# Applies the workbook update for the newly-closed trade (global trade #188,
# "momentum" strategy) plus the two freshly-opened trades (#217 on
# HighProbConvergence, #218 on MarketMaking) to every sheet that tracks them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet - aggregate counters
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 188     # Total Trades
$summary.Range("B9").Value = 44.15   # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - "momentum" row (row 11)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D11").Value = 49      # Trades
$status.Range("G11").Value = 30.61   # Win rate %

# ---------------------------------------------------------------------
# All Trades sheet - close out trade #188 (row 189) and append the two
# newly-opened trades (#217 / #218) as rows 218 / 219
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Cells.Item(189, 7).Value = 0.01          # G189 Exit Price
$allTrades.Cells.Item(189, 8).Value = "CLOSED"      # H189 Status
$allTrades.Cells.Item(189, 11).Value = 99.22        # K189 Capital After
$allTrades.Cells.Item(189, 12).Value = "early_exit" # L189 Exit Reason
$allTrades.Cells.Item(189, 13).Value = 0.17         # M189 Duration (min)

# New row 218 -> trade #217 (HighProbConvergence, UP)
$allTrades.Cells.Item(218, 1).Value = 217
$allTrades.Cells.Item(218, 2).Value = "'2026-02-18"
$allTrades.Cells.Item(218, 3).Value = "00:44:10"
$allTrades.Cells.Item(218, 4).Value = "HighProbConvergence"
$allTrades.Cells.Item(218, 5).Value = "UP"
$allTrades.Cells.Item(218, 6).Value = 0.01
$allTrades.Cells.Item(218, 8).Value = "OPEN"
$allTrades.Cells.Item(218, 9).Value = 0
$allTrades.Cells.Item(218, 10).Value = 0
$allTrades.Cells.Item(218, 11).Value = 100.3223499536821
$allTrades.Cells.Item(218, 13).Value = 0
$allTrades.Cells.Item(218, 14).Value = 0
$allTrades.Cells.Item(218, 15).Value = 0
$allTrades.Cells.Item(218, 16).Value = 0.95
$allTrades.Cells.Item(218, 17).Value = "Mean reversion UP: price 11.35% below mean (z=-4.36)"

# New row 219 -> trade #218 (MarketMaking, UP)
$allTrades.Cells.Item(219, 1).Value = 218
$allTrades.Cells.Item(219, 2).Value = "'2026-02-18"
$allTrades.Cells.Item(219, 3).Value = "00:44:11"
$allTrades.Cells.Item(219, 4).Value = "MarketMaking"
$allTrades.Cells.Item(219, 5).Value = "UP"
$allTrades.Cells.Item(219, 6).Value = 0.01
$allTrades.Cells.Item(219, 8).Value = "OPEN"
$allTrades.Cells.Item(219, 9).Value = 0
$allTrades.Cells.Item(219, 10).Value = 0
$allTrades.Cells.Item(219, 11).Value = 99.45858346467946
$allTrades.Cells.Item(219, 13).Value = 0
$allTrades.Cells.Item(219, 14).Value = 0
$allTrades.Cells.Item(219, 15).Value = 0
$allTrades.Cells.Item(219, 16).Value = 0.6
$allTrades.Cells.Item(219, 17).Value = "Normal spread capture: 225 bps"

# ---------------------------------------------------------------------
# momentum sheet - same trade #188 close-out lives in row 50 here, but the
# per-strategy sheets use a different column layout:
# L=Entry Slippage, M=Exit Slippage, N=Confidence, O=Entry Reason,
# P=Exit Reason, Q=Duration (min)
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Cells.Item(50, 7).Value = 0.01           # G50 Exit Price
$momentum.Cells.Item(50, 8).Value = "CLOSED"       # H50 Status
$momentum.Cells.Item(50, 11).Value = 99.22         # K50 Capital After
$momentum.Cells.Item(50, 16).Value = "early_exit"  # P50 Exit Reason
$momentum.Cells.Item(50, 17).Value = 0.17          # Q50 Duration (min)

# ---------------------------------------------------------------------
# HighProbConvergence sheet - append newly-opened trade #217 as row 27
# ---------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")
$hpc.Cells.Item(27, 1).Value = 217
$hpc.Cells.Item(27, 2).Value = "'2026-02-18"
$hpc.Cells.Item(27, 3).Value = "00:44:10"
$hpc.Cells.Item(27, 4).Value = "HighProbConvergence"
$hpc.Cells.Item(27, 5).Value = "UP"
$hpc.Cells.Item(27, 6).Value = 0.01
$hpc.Cells.Item(27, 8).Value = "OPEN"
$hpc.Cells.Item(27, 9).Value = 0
$hpc.Cells.Item(27, 10).Value = 0
$hpc.Cells.Item(27, 11).Value = 100.3223499536821
$hpc.Cells.Item(27, 12).Value = 0
$hpc.Cells.Item(27, 13).Value = 0
$hpc.Cells.Item(27, 14).Value = 0.95
$hpc.Cells.Item(27, 15).Value = "Mean reversion UP: price 11.35% below mean (z=-4.36)"
$hpc.Cells.Item(27, 17).Value = 0

# ---------------------------------------------------------------------
# MarketMaking sheet - append newly-opened trade #218 as row 102
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Cells.Item(102, 1).Value = 218
$mm.Cells.Item(102, 2).Value = "'2026-02-18"
$mm.Cells.Item(102, 3).Value = "00:44:11"
$mm.Cells.Item(102, 4).Value = "MarketMaking"
$mm.Cells.Item(102, 5).Value = "UP"
$mm.Cells.Item(102, 6).Value = 0.01
$mm.Cells.Item(102, 8).Value = "OPEN"
$mm.Cells.Item(102, 9).Value = 0
$mm.Cells.Item(102, 10).Value = 0
$mm.Cells.Item(102, 11).Value = 99.45858346467946
$mm.Cells.Item(102, 12).Value = 0
$mm.Cells.Item(102, 13).Value = 0
$mm.Cells.Item(102, 14).Value = 0.6
$mm.Cells.Item(102, 15).Value = "Normal spread capture: 225 bps"
$mm.Cells.Item(102, 17).Value = 0
